$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H69").Value = 4799.2

$ws.Range("H72").Value = 4799.2

$ws.Range("H138").Value = 5734.2607
$ws.Range("I138").Value = 5048.5
$ws.Range("J138").Value = 5878.6313
$ws.Range("K138").Value = 15145.5
$ws.Range("L138").Value = 17635.8939
$ws.Range("M138").Value = -10005.5
$ws.Range("N138").Value = -27915.8939

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2060.625
$ws.Range("I2").Value = 2051.5386
$ws.Range("J2").Value = 2100
$ws.Range("K2").Value = 2051.5386
$ws.Range("L2").Value = 2100
$ws.Range("M2").Value = -1938.5386
$ws.Range("N2").Value = -2326

$ws.Range("H61").Value = 3831.3333
$ws.Range("I61").Value = 1495
$ws.Range("K61").Value = 1495
$ws.Range("M61").Value = -1283

$ws.Range("H116").Value = 2060.625
$ws.Range("I116").Value = 2051.5386
$ws.Range("J116").Value = 2100
$ws.Range("K116").Value = 2051.5386
$ws.Range("L116").Value = 2100
$ws.Range("M116").Value = 242.4614000000001
$ws.Range("N116").Value = -6688

$ws.Range("H136").Value = 3831.3333
$ws.Range("I136").Value = 1495
$ws.Range("K136").Value = 4485
$ws.Range("M136").Value = -1935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2060.625
$ws.Range("I3").Value = 2051.5386
$ws.Range("J3").Value = 2100
$ws.Range("K3").Value = 2051.5386
$ws.Range("L3").Value = 2100
$ws.Range("M3").Value = -1937.5386
$ws.Range("N3").Value = -2328

$ws.Range("H134").Value = 935
$ws.Range("I134").Value = 935
$ws.Range("K134").Value = 2805
$ws.Range("M134").Value = -270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 27857.143
$ws.Range("J44").Value = 27857.143
$ws.Range("L44").Value = 27857.143
$ws.Range("N44").Value = -28741.143

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H55").Value = 5000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 5000
$ws.Range("N55").Value = -5630
$ws.Range("M55").ClearContents()

$ws.Range("H58").Value = 5703.3335
$ws.Range("I58").Value = 7555.5
$ws.Range("J58").Value = 1999
$ws.Range("K58").Value = 7555.5
$ws.Range("L58").Value = 1999
$ws.Range("M58").Value = -7352.5
$ws.Range("N58").Value = -2405

$ws.Range("H99").Value = 3004.7778
$ws.Range("I99").Value = 2674.1667
$ws.Range("J99").Value = 3666
$ws.Range("K99").Value = 2674.1667
$ws.Range("L99").Value = 3666
$ws.Range("M99").Value = -1176.1667
$ws.Range("N99").Value = -6662

$ws.Range("H126").Value = 3004.7778
$ws.Range("I126").Value = 2674.1667
$ws.Range("J126").Value = 3666
$ws.Range("K126").Value = 8022.500100000001
$ws.Range("L126").Value = 10998
$ws.Range("M126").Value = -5552.500100000001
$ws.Range("N126").Value = -15938

$ws.Range("H132").Value = 3227.818
$ws.Range("I132").Value = 3300.7
$ws.Range("K132").Value = 9902.099999999999
$ws.Range("M132").Value = -7372.099999999999

$ws.Range("H136").Value = 5703.3335
$ws.Range("I136").Value = 7555.5
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 22666.5
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -20116.5
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 974.6842
$ws.Range("J131").Value = 990
$ws.Range("L131").Value = 2970
$ws.Range("N131").Value = -13050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3316
$ws.Range("I41").Value = 3579.2
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 3579.2
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = -3224.2
$ws.Range("N41").Value = -2710

$ws.Range("H46").Value = 7460.4443
$ws.Range("I46").Value = 2716.6667
$ws.Range("J46").Value = 9832.333000000001
$ws.Range("K46").Value = 2716.6667
$ws.Range("L46").Value = 9832.333000000001
$ws.Range("M46").Value = -2560.6667
$ws.Range("N46").Value = -10144.333

$ws.Range("H122").Value = 3744.111
$ws.Range("I122").Value = 3587.125
$ws.Range("K122").Value = 10761.375
$ws.Range("M122").Value = -8311.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H82").Value = 1605.8889
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 1681.625
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 1681.625
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -2403.625

$ws.Range("H85").Value = 1605.8889
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 1681.625
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 1681.625
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -4177.625

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H103").Value = 23244.666
$ws.Range("J103").Value = 23244.666
$ws.Range("L103").Value = 23244.666
$ws.Range("N103").Value = -25588.666

$ws.Range("H122").Value = 6684
$ws.Range("I122").Value = 5995.3335
$ws.Range("K122").Value = 17986.0005
$ws.Range("M122").Value = -15536.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2530.6667
$ws.Range("I81").Value = 2530.6667
$ws.Range("K81").Value = 5061.3334
$ws.Range("M81").Value = -4000.3334

$ws.Range("H84").Value = 2530.6667
$ws.Range("I84").Value = 2530.6667
$ws.Range("K84").Value = 25306.667
$ws.Range("M84").Value = -20002.667

$ws.Range("H97").Value = 34142.25
$ws.Range("J97").Value = 34142.25
$ws.Range("L97").Value = 34142.25
$ws.Range("N97").Value = -36124.25

$ws.Range("H107").Value = 1233.0667
$ws.Range("I107").Value = 1081
$ws.Range("J107").Value = 1309.1
$ws.Range("K107").Value = 3243
$ws.Range("L107").Value = 3927.3
$ws.Range("M107").Value = -1323
$ws.Range("N107").Value = -7767.299999999999

$ws.Range("H122").Value = 2656.1428
$ws.Range("I122").Value = 2898.8333
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 8696.499899999999
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -6246.499899999999
$ws.Range("N122").Value = -8500

$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140
